$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 5-17, keeping only header (row1) and rows 2-4
$ws.Range("A5:B17").EntireRow.Delete() | Out-Null

# Update row 2
$ws.Range("A2").Value = "05/13/2021 20:18:13"
$ws.Range("B2").Value = 30.642

# Update row 3
$ws.Range("A3").Value = "05/13/2021 20:23:32"
$ws.Range("B3").Value = 13.236

# Update row 4
$ws.Range("A4").Value = "05/13/2021 20:23:51"
$ws.Range("B4").Value = 8.685
